# Back_Testing_Results_Sample.xlsx -- "Results in 920 revamped"
#
# Summary of changes:
#  - Trades sheet gets 5 new header columns (F:J) with new labels
#  - Trades sheet becomes the active/selected sheet, with its own view state
#  - Summary sheet view scrolls down (topLeftCell) and selection moves to B13:K23
#
$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("Summary")
$wsTrades  = $wb.Worksheets.Item("Trades")

# --- Trades sheet: add the new header cells (columns F..J on row 1) ---
# Order matters for shared-string table layout: F, G, I, J, H
$wsTrades.Range("F1").Value = "Profit/Loss"
$wsTrades.Range("G1").Value = "Reverse Entry"
$wsTrades.Range("I1").Value = "Reverse_Exit"
$wsTrades.Range("J1").Value = "R_Close_Exit"
$wsTrades.Range("H1").Value = "R_Close_Entry"

# --- Trades sheet: column widths for the new layout ---
# (ColumnWidth is quantized by the host to a 1/6-character grid, so these are
#  the inputs whose rounded/stored width lands closest to the authored widths.)
$wsTrades.Columns.Item(2).ColumnWidth = 14.833333333333334
$wsTrades.Columns.Item(4).ColumnWidth = 9.333333333333334
$wsTrades.Columns.Item(6).ColumnWidth = 15.833333333333334

# --- Trades sheet: view state (becomes the active sheet) ---
$wsTrades.Activate()
$wsTrades.Range("H2").Select()

# --- Summary sheet: view state (scrolled down, new selection) ---
$wsSummary.Activate()
$excel.ActiveWindow.ScrollRow = 8
$wsSummary.Range("B13:K23").Select()

# --- Make Trades the active tab, as in the target workbook ---
$wsTrades.Activate()
